$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I and J ---
# Copy formatting (bold font + border + alignment) from an existing header
# cell (H1) onto the new header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-20: new numeric columns I (I0) and J (IF) ---
$data = @{
    2  = @(8, 9)
    3  = @(6, 7)
    4  = @(6, 7)
    5  = @(8, 8)
    6  = @(8, 9)
    7  = @(9, 9)
    8  = @(6, 6)
    9  = @(7, 9)
    10 = @(6, 8)
    11 = @(9, 9)
    12 = @(2, 4)
    13 = @(7, 8)
    14 = @(5, 7)
    15 = @(7, 8)
    16 = @(6, 9)
    17 = @(8, 8)
    18 = @(8, 9)
    19 = @(8, 9)
    20 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Output "done"
